$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1377.2727
$ws.Range("I98").Value = 1166.6666
$ws.Range("J98").Value = 1630
$ws.Range("K98").Value = 1166.6666
$ws.Range("L98").Value = 1630
$ws.Range("M98").Value = 331.3334
$ws.Range("N98").Value = -4626

$ws.Range("H122").Value = 1377.2727
$ws.Range("I122").Value = 1166.6666
$ws.Range("J122").Value = 1630
$ws.Range("K122").Value = 3499.9998
$ws.Range("L122").Value = 4890
$ws.Range("M122").Value = -1049.9998
$ws.Range("N122").Value = -9790

$ws.Range("H132").Value = 3368.1333
$ws.Range("I132").Value = 3184
$ws.Range("J132").Value = 4565
$ws.Range("K132").Value = 9552
$ws.Range("L132").Value = 13695
$ws.Range("M132").Value = -7022
$ws.Range("N132").Value = -18755

$ws.Range("H137").Value = 3334575.8
$ws.Range("I137").Value = 1516352
$ws.Range("J137").Value = 8334691.5
$ws.Range("K137").Value = 4549056
$ws.Range("L137").Value = 25004074.5
$ws.Range("M137").Value = -4546506
$ws.Range("N137").Value = -25009174.5

$ws.Range("H138").Value = 1907.7894
$ws.Range("I138").Value = 1205.7059
$ws.Range("J138").Value = 2945.652
$ws.Range("K138").Value = 3617.1177
$ws.Range("L138").Value = 8836.956
$ws.Range("M138").Value = 1522.8823
$ws.Range("N138").Value = -19116.956

$ws.Range("H141").Value = 1614.8414
$ws.Range("I141").Value = 1098.0227
$ws.Range("J141").Value = 2213.2632
$ws.Range("K141").Value = 3294.0681
$ws.Range("L141").Value = 6639.7896
$ws.Range("M141").Value = 1885.9319
$ws.Range("N141").Value = -16999.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3434159
$ws.Range("I32").Value = 5718.476
$ws.Range("J32").Value = 25033334
$ws.Range("K32").Value = 5718.476
$ws.Range("L32").Value = 25033334
$ws.Range("M32").Value = -5431.476
$ws.Range("N32").Value = -25033908

$ws.Range("H45").Value = 2886
$ws.Range("I45").Value = 1205.2667
$ws.Range("J45").Value = 6487.5713
$ws.Range("K45").Value = 1205.2667
$ws.Range("L45").Value = 6487.5713
$ws.Range("M45").Value = -828.2666999999999
$ws.Range("N45").Value = -7241.5713

$ws.Range("H52").Value = 16162.5
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 16162.5
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 16162.5
$ws.Range("N52").Value = -16798.5

$ws.Range("H122").Value = 1884.3214
$ws.Range("I122").Value = 1780.6842
$ws.Range("J122").Value = 2103.111
$ws.Range("K122").Value = 5342.0526
$ws.Range("L122").Value = 6309.333
$ws.Range("M122").Value = -2892.0526
$ws.Range("N122").Value = -11209.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 36000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 36000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 36000
$ws.Range("N18").Value = -36460

$ws.Range("H58").Value = 2430.7144
$ws.Range("I58").Value = 2263
$ws.Range("J58").Value = 2850
$ws.Range("K58").Value = 2263
$ws.Range("L58").Value = 2850
$ws.Range("M58").Value = -2060
$ws.Range("N58").Value = -3256

$ws.Range("H132").Value = 2941.2927
$ws.Range("I132").Value = 2385.889
$ws.Range("J132").Value = 6940.2
$ws.Range("K132").Value = 7157.667
$ws.Range("L132").Value = 20820.6
$ws.Range("M132").Value = -4627.667
$ws.Range("N132").Value = -25880.6

$ws.Range("H134").Value = 4894.4287
$ws.Range("I134").Value = 5555.237
$ws.Range("J134").Value = 2611.6365
$ws.Range("K134").Value = 16665.711
$ws.Range("L134").Value = 7834.9095
$ws.Range("M134").Value = -14130.711
$ws.Range("N134").Value = -12904.9095

$ws.Range("H136").Value = 2430.7144
$ws.Range("I136").Value = 2263
$ws.Range("J136").Value = 2850
$ws.Range("K136").Value = 6789
$ws.Range("L136").Value = 8550
$ws.Range("M136").Value = -4239
$ws.Range("N136").Value = -13650

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 916.6667
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 1060
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 3180
$ws.Range("M22").Value = -431
$ws.Range("N22").Value = -3518

$ws.Range("H27").Value = 916.6667
$ws.Range("I27").Value = 200
$ws.Range("J27").Value = 1060
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 3180
$ws.Range("M27").Value = -498
$ws.Range("N27").Value = -3384

$ws.Range("H34").Value = 1146.8
$ws.Range("I34").Value = 206.66667
$ws.Range("J34").Value = 1549.7142
$ws.Range("K34").Value = 620.00001
$ws.Range("L34").Value = 4649.142599999999
$ws.Range("M34").Value = -536.00001
$ws.Range("N34").Value = -4817.142599999999

$ws.Range("H40").Value = 147.71428
$ws.Range("I40").Value = 159
$ws.Range("J40").Value = 80
$ws.Range("K40").Value = 636
$ws.Range("L40").Value = 320
$ws.Range("M40").Value = -567
$ws.Range("N40").Value = -458

$ws.Range("H58").Value = 4204092.5
$ws.Range("I58").Value = 7354537

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3512
$ws.Range("I132").Value = 3512
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10536
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8006
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1673.6444
$ws.Range("I136").Value = 1660.5428
$ws.Range("J136").Value = 1719.5
$ws.Range("K136").Value = 4981.6284
$ws.Range("L136").Value = 5158.5
$ws.Range("M136").Value = -2431.6284
$ws.Range("N136").Value = -10258.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 11000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 11000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 11000
$ws.Range("N19").Value = -11348

$ws.Range("H40").Value = 18000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 18000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 18000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -18298

$ws.Range("H57").Value = 14999.667
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 14999.667
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 14999.667
$ws.Range("N57").Value = -16507.667

$ws.Range("H62").Value = 3621.111
$ws.Range("I62").Value = 3570
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 3570
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -2946
$ws.Range("N62").Value = -5048

$ws.Range("H65").Value = 3621.111
$ws.Range("I65").Value = 3570
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 17850
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -14730
$ws.Range("N65").Value = -25240

$ws.Range("H108").Value = 10300
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 10300
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 10300
$ws.Range("N108").Value = -17980

$ws.Range("H132").Value = 1677.1666
$ws.Range("I132").Value = 2111.5
$ws.Range("J132").Value = 808.5
$ws.Range("K132").Value = 6334.5
$ws.Range("L132").Value = 2425.5
$ws.Range("M132").Value = -3804.5
$ws.Range("N132").Value = -7485.5

$ws.Range("H136").Value = 2215.4146
$ws.Range("I136").Value = 2104.5293
$ws.Range("J136").Value = 2754
$ws.Range("K136").Value = 6313.5879
$ws.Range("L136").Value = 8262
$ws.Range("M136").Value = -3763.5879
$ws.Range("N136").Value = -13362
